$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: praclen 4 -> 5
$ws.Cells.Item(2,9).Value = 5

# Row 3: praclen 4 -> 5
$ws.Cells.Item(3,9).Value = 5

# Row 4: updated schedule values
$ws.Cells.Item(4,2).Value = 3
$ws.Cells.Item(4,3).Value = 2
$ws.Cells.Item(4,4).Value = 3
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 61
$ws.Cells.Item(4,9).Value = 5

# Row 5: updated schedule values
$ws.Cells.Item(5,2).Value = 2
$ws.Cells.Item(5,3).Value = 3
$ws.Cells.Item(5,4).Value = 2
$ws.Cells.Item(5,7).Value = 3
$ws.Cells.Item(5,8).Value = 64
$ws.Cells.Item(5,9).Value = 5

# New row 6, matching the old row 5 values (praclen updated to 5)
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = 1
$ws.Cells.Item(6,4).Value = 1
$ws.Cells.Item(6,5).Value = 6
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 5
$ws.Cells.Item(6,8).Value = 66
$ws.Cells.Item(6,9).Value = 5
$ws.Cells.Item(6,10).Value = "train_dim1_1"

# Update selection to match final state
$ws.Range("L14").Select()
